{"js": "const body = context.document.body;\n\n// 1. The last bullet item currently ends with a period (\"...\u043e\u0442\u0432\u0435\u0440\u0441\u0442\u0438\u044f.\").\n//    Change it to a semicolon so it matches the punctuation style used by\n//    the rest of the list.\nconst results = body.search(\"\u043e\u0442\u0432\u0435\u0440\u0441\u0442\u0438\u044f.\", { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"\u043e\u0442\u0432\u0435\u0440\u0441\u0442\u0438\u044f;\", Word.InsertLocation.replace);\n}\n\n// 2. Drop the stale \"_GoBack\" bookmark that sat at the end of that paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 3. Append a new bulleted list item with the additional task.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\n  \"\u0412\u044b\u0440\u043e\u0432\u043d\u044f\u0442\u044c \u0434\u043b\u0438\u043d\u044b \u0432\u0441\u0435\u0445 \u0434\u043e\u0440\u043e\u0436\u0435\u043a (\u0445\u043e\u0442\u044f \u0431\u044b \u0428\u0418\u041c-\u043a\u0430\u043d\u0430\u043b\u044b).\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The last bullet item currently ends with a period (\"...\u043e\u0442\u0432\u0435\u0440\u0441\u0442\u0438\u044f.\") and\n# carries a leftover \"_GoBack\" bookmark. Change the trailing period to a\n# semicolon (to match the other list items) and drop the bookmark.\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$find = $lastPara.Range.Find\n$find.Execute(\"\u043e\u0442\u0432\u0435\u0440\u0441\u0442\u0438\u044f.\", $false, $false, $false, $false, $false, $true, 0, $false, \"\u043e\u0442\u0432\u0435\u0440\u0441\u0442\u0438\u044f;\", 2)\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Append a brand-new bulleted paragraph (same list) with the new task.\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$lastPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs($d.Paragraphs.Count)\n$newPara.Range.Text = \"\u0412\u044b\u0440\u043e\u0432\u043d\u044f\u0442\u044c \u0434\u043b\u0438\u043d\u044b \u0432\u0441\u0435\u0445 \u0434\u043e\u0440\u043e\u0436\u0435\u043a (\u0445\u043e\u0442\u044f \u0431\u044b \u0428\u0418\u041c-\u043a\u0430\u043d\u0430\u043b\u044b).\"\n"}
